# Add a "Sheet2" right after Sheet1, populate it with a copy of Sheet1's
# header row + first data row (CAGR / "df" columns included), but with the
# "variable" value changed to "z", then leave Sheet2 as the active sheet
# with A2 selected, and leave Sheet1 with A1:O2 selected.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New sheet, placed immediately after Sheet1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# ---- Row 1: header row (identical to Sheet1's header row) ----
$ws2.Range("A1").Value = "variable"
$ws2.Range("B1").Value = "scenario"
$ws2.Range("C1").Value = "module"
$ws2.Range("D1").Value = "distribution"
$ws2.Range("E1").Value = "param 1"
$ws2.Range("F1").Value = "param 2"

$ws2.Range("G1").Value = "param 3"
$ws2.Range("H1").Value = "unit"
$ws2.Range("I1").Value = "start date"
$ws2.Range("J1").Value = "end date"
$ws2.Range("K1").Value = "CAGR"
$ws2.Range("L1").Value = "ref date"
$ws2.Range("G1:L1").Font.Color = 0

$ws2.Range("M1").Value = "label"

$ws2.Range("N1").Value = "comment"
$ws2.Range("O1").Value = "source"
$ws2.Range("N1:O1").Font.Color = 0

# ---- Row 2: first data row (same as Sheet1's row 2, but variable = "z") ----
$ws2.Range("A2").Value = "z"
$ws2.Range("C2").Value = "numpy.random"
$ws2.Range("D2").Value = "choice"
$ws2.Range("E2").Value = 1

$ws2.Range("H2").Value = "kg"
$ws2.Range("H2").Font.Color = 0

$ws2.Range("I2").NumberFormat = "m/d/yy"
$ws2.Range("I2").Value = 39814
$ws2.Range("J2").NumberFormat = "m/d/yy"
$ws2.Range("J2").Value = 39904
$ws2.Range("K2").NumberFormat = "0.00"
$ws2.Range("K2").Value = 0.1
$ws2.Range("L2").NumberFormat = "m/d/yy"
$ws2.Range("L2").Value = 39814

$ws2.Range("M2").Value = "test var 1"

# ---- Selections ----
# Sheet1 keeps a block selection covering the header + first data row.
$ws1.Range("A1:O2").Select()

# Sheet2 is the sheet the user ends up looking at, with A2 selected.
$ws2.Activate()
$ws2.Range("A2").Select()
